# Update all example folders to newest format / Switch to newest InOutModule
#
# The new workbook format gives the title row (row 1) on every data sheet
# a taller, explicit row height (24 points) instead of the default height,
# so the bigger title font has room to breathe.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).RowHeight = 24
}
